$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 9, pushing the existing data rows
# (old rows 9-28) down to rows 11-30.
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(10).Insert()

# Fill in the new row 9 with this week's first new data point.
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 44497
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 100114002
$ws.Range("G9").Value = "Camote"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 20000
$ws.Range("N9").Value = "$/caja 15 kilos granel"
$ws.Range("O9").Value = "Perú"
$ws.Range("P9").Value = 1333
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"

# Fill in the new row 10 with this week's second new data point.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44497
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100114002
$ws.Range("G10").Value = "Camote"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("N10").Value = "$/malla 20 kilos"
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 1000
$ws.Range("Q10").Value = 20
$ws.Range("R10").Value = "Hortaliza"
